$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 98; $r++) {
    $ws.Range("K$r").Formula = "=E$r/D$r"
    $ws.Range("L$r").Formula = "=H$r/F$r"
}

[void]$ws.Range("O6").Select()
